# Applies the "first version that should work on server" update to the
# Work-Report Trello sheet: refreshed counters, a round of name swaps
# between Theresa Rinnert / Theresa Schmid, several newly-listed helpers,
# and renamed "Karten ohne Aktivität" entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlPasteValues = -4163

# --- Aktiv block (rows 17-21) ---
$ws.Range("C17").Value = 15
$ws.Range("G17").Value = 12

# C18 keeps holding a textual "12" (not a real number) in the source report,
# so force it back to text after the write instead of letting it settle as
# a number.
$ws.Range("C18").Formula = '="12"'
$ws.Range("C18").Copy()
$ws.Range("C18").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = $false

$ws.Range("G18").Value = 2

# Same trick for D19, which also stores a plain count as text.
$ws.Range("D19").Formula = '="4"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = $false

$ws.Range("G19").Value = 1
$ws.Range("D20").Value = "80.0% der Karten"
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = 3

# --- Karten ohne Aktivität block (rows 27-31) ---
$ws.Range("C27").Value = 36

$ws.Range("B28").Value = "Informationen beschaffen für den Durchlauf mit dem Gastronomiezweig 💩"
$ws.Range("C28").Value = 26

$ws.Range("B29").Value = "BB 2.0 - Kriterien für Bewertung/Auswahl von Partnern 💩"
$ws.Range("C29").Value = 26

$ws.Range("B30").Value = "BB 2.0 - Risikobewertung 💩"
$ws.Range("C30").Value = 26

$ws.Range("B31").Value = "BB 2.0 - Kostenplan 💩"
$ws.Range("C31").Value = 26

# --- Gemeinschaftlich / Aktivste Mitglieder block ---
$ws.Range("F34").Value = 6
$ws.Range("G34").Value = "(40.0%)"

$ws.Range("B36").Value = "Theresa Schmid"
$ws.Range("B37").Value = "Theresa Rinnert"
$ws.Range("C37").Value = 3

# Row 39 previously had no entry in column B; add the new helper and copy
# the (unstyled) format from a sibling cell in the same list so no new
# style gets invented for it.
$ws.Range("B39").Value = "Nina Stallmann"
$ws.Range("B70").Copy()
$ws.Range("B39").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("C39").Value = 1

$ws.Range("E40").Value = "Theresa Schmid"
$ws.Range("E41").Value = "Theresa Rinnert"
$ws.Range("F41").Value = 3

$ws.Range("E43").Value = "Nina Stallmann"
$ws.Range("E40").Copy()
$ws.Range("E43").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F43").Value = 1

# --- Zuverlässig / Aktivste Helfer block ---
$ws.Range("B48").Value = "Marie-Sophie Braun"
$ws.Range("B70").Copy()
$ws.Range("B48").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("C48").Value = 1

$ws.Range("E50").Value = "Jacob Escherich"
$ws.Range("F50").Value = 1
$ws.Range("F34").Copy()
$ws.Range("F50").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("E51").Value = "Micha Landoll"
$ws.Range("E40").Copy()
$ws.Range("E51").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F51").Value = 1

# --- Pünktlich block ---
$ws.Range("A60").Value = "Marie-Sophie Braun"
$ws.Range("B60").Value = 1
$ws.Range("G60").Value = 1

$ws.Range("G61").Value = 1

# --- Detailliert block ---
$ws.Range("B70").Value = 15
$ws.Range("B73").Value = 8
